$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Exam 2 weighting update: the running total (E4) now only counts half
# of the raw point sum (mid-semester re-weighting).
$ws.Range("E4").Formula = "=SUM(B2:B30)*0.5"

# Fix the "out of" label text for the 1350-point track.
$ws.Range("J10").Value = "Out of 1350"

# New "Total Points (To Date)" row.
$ws.Range("D9").Value = "Total Points"
$ws.Range("E9").Formula = "=E4"
$ws.Range("D10").Value = "(To Date)"

# Leave the cursor where the author left it.
$ws.Range("E10").Select()
